function Set-TextCell($ws, $addr, $val) {
    # Force the value to be stored as text even when it looks like a
    # number (e.g. "255.12"), matching the source data which keeps
    # these price columns as plain strings, not numeric cells.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "D2" '99.107.84'
$ws.Range("E2").Value = '  +1.66%  '
Set-TextCell $ws "D3" '3.317.00'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextCell $ws "D5" '255.12'
$ws.Range("E5").Value = '  -0.01%  '
Set-TextCell $ws "D6" '624.76'
$ws.Range("E6").Value = '  +0.89%  '
Set-TextCell $ws "D7" '1.46'
$ws.Range("E7").Value = '  +31.13%  '
Set-TextCell $ws "D8" '0.411'
$ws.Range("E8").Value = '  +7.09%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  +22.24%  '
Set-TextCell $ws "D11" '3.315.25'
$ws.Range("E11").Value = '  -0.35%  '
Set-TextCell $ws "D12" '0.201'
$ws.Range("E12").Value = '  +0.82%  '
Set-TextCell $ws "D13" '39.32'
$ws.Range("E13").Value = '  +11.11%  '
Set-TextCell $ws "D14" '98.720.94'
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("E15").Value = '  +2.42%  '
Set-TextCell $ws "D16" '3.939.16'
$ws.Range("E16").Value = '  +0.57%  '
Set-TextCell $ws "D17" '5.49'
$ws.Range("E17").Value = '  -0.40%  '
Set-TextCell $ws "D18" '3.316.64'
$ws.Range("E18").Value = '  -0.19%  '
Set-TextCell $ws "D19" '3.49'
$ws.Range("E19").Value = '  -2.00%  '
Set-TextCell $ws "D20" '15.61'
$ws.Range("E20").Value = '  +4.36%  '
Set-TextCell $ws "D21" '6.31'
$ws.Range("E21").Value = '  +8.91%  '
Set-TextCell $ws "D22" '487.92'
$ws.Range("E22").Value = '  +1.36%  '
Set-TextCell $ws "D23" '9.49'
$ws.Range("E23").Value = '  +3.04%  '
Set-TextCell $ws "D24" '0.0000204'
$ws.Range("E24").Value = '  -1.95%  '
Set-TextCell $ws "D25" '5.65'
$ws.Range("E25").Value = '  +0.49%  '
Set-TextCell $ws "D26" '88.90'
$ws.Range("E26").Value = '  +1.45%  '
Set-TextCell $ws "D27" '12.03'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("E28").Value = '  +26.98%  '
Set-TextCell $ws "D29" '3.488.55'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws "D31" '0.189'
$ws.Range("E31").Value = '  +3.41%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws "D32" '0.136'
$ws.Range("E32").Value = '  +11.46%  '
Set-TextCell $ws "D33" '10.13'
$ws.Range("E33").Value = '  +10.03%  '
$ws.Range("E34").Value = '  +0.18%  '
Set-TextCell $ws "D35" '27.89'
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws "D36" '0.149'
$ws.Range("E36").Value = '  -1.23%  '
Set-TextCell $ws "D37" '7.24'
$ws.Range("E37").Value = '  -2.14%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextCell $ws "D38" '0.470'
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("E40").Value = '  +0.18%  '
Set-TextCell $ws "D41" '492.43'
$ws.Range("E41").Value = '  -2.96%  '
Set-TextCell $ws "D42" '3.65'
$ws.Range("E42").Value = '  +3.33%  '
$ws.Range("E43").Value = '  -2.75%  '
Set-TextCell $ws "D44" '0.792'
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("E45").Value = '  +0.03%  '
Set-TextCell $ws "D46" '3.14'
$ws.Range("E46").Value = '  -4.50%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws "D47" '1.98'
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws "D48" '159.46'
$ws.Range("E48").Value = '  -0.29%  '
Set-TextCell $ws "D49" '7.34'
$ws.Range("E49").Value = '  +16.63%  '
Set-TextCell $ws "D50" '0.850'
$ws.Range("E50").Value = '  +6.86%  '
Set-TextCell $ws "D51" '4.74'
$ws.Range("E51").Value = '  +5.35%  '
